# Address code review feedback - improve compatibility and maintainability
# 1) Rename the summary sheet so the English label leads (drop the emoji prefix)
# 2) Bump the "created at" timestamp
# 3) Replace the emoji-coded status labels in column E with plain bilingual text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet tab: "📊 الملخص - Summary" -> "Summary - الملخص"
$ws.Name = "Summary - الملخص"

# Update the creation timestamp in A4
$ws.Range("A4").Value = "تاريخ الإنشاء: 2025-12-09 06:56"

# Map each status row to its new plain-text label (emoji removed, English suffix added)
$statusByRow = @{
    7  = "كبير جداً - Very Large"
    8  = "كبير - Large"
    9  = "كبير - Large"
    10 = "كبير - Large"
    11 = "كبير - Large"
    12 = "كبير - Large"
    13 = "كبير - Large"
    14 = "كبير - Large"
    15 = "كبير - Large"
    16 = "متوسط - Medium"
    17 = "متوسط - Medium"
    18 = "متوسط - Medium"
    19 = "متوسط - Medium"
    20 = "متوسط - Medium"
    21 = "صغير - Small"
    22 = "صغير - Small"
    23 = "صغير - Small"
    24 = "صغير - Small"
    25 = "صغير - Small"
    26 = "صغير - Small"
    27 = "صغير - Small"
    28 = "صغير - Small"
    29 = "صغير - Small"
    30 = "صغير - Small"
    31 = "صغير - Small"
    32 = "صغير - Small"
}

foreach ($row in $statusByRow.Keys) {
    $ws.Cells.Item($row, 5).Value = $statusByRow[$row]
}
